{"js": "// Replace the three-digit x one-digit multiplication prompts throughout\n// the document body (table cells) with their new values.\nconst replacements = [\n  [\"530\u00d78=\", \"740\u00d74=\"],\n  [\"249\u00d74=\", \"733\u00d74=\"],\n  [\"477\u00d76=\", \"853\u00d77=\"],\n  [\"788\u00d75=\", \"871\u00d79=\"],\n  [\"988\u00d75=\", \"321\u00d74=\"],\n  [\"546\u00d76=\", \"553\u00d73=\"],\n  [\"120\u00d73=\", \"803\u00d75=\"],\n  [\"551\u00d72=\", \"965\u00d76=\"],\n  [\"634\u00d75=\", \"361\u00d73=\"],\n  [\"996\u00d75=\", \"126\u00d78=\"],\n  [\"648\u00d79=\", \"996\u00d74=\"],\n  [\"631\u00d74=\", \"418\u00d74=\"],\n  [\"597\u00d77=\", \"710\u00d76=\"],\n  [\"974\u00d77=\", \"166\u00d79=\"],\n  [\"580\u00d78=\", \"486\u00d74=\"],\n  [\"175\u00d72=\", \"654\u00d75=\"],\n  [\"250\u00d76=\", \"538\u00d76=\"],\n  [\"576\u00d78=\", \"486\u00d79=\"],\n  [\"141\u00d73=\", \"695\u00d75=\"],\n  [\"390\u00d75=\", \"147\u00d76=\"],\n  [\"863\u00d78=\", \"611\u00d73=\"],\n  [\"672\u00d72=\", \"678\u00d72=\"],\n  [\"813\u00d75=\", \"435\u00d76=\"],\n  [\"415\u00d78=\", \"483\u00d77=\"],\n  [\"378\u00d77=\", \"186\u00d73=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (const r of results.items) {\n    r.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"530\u00d78=\"; New = \"740\u00d74=\" },\n    @{ Old = \"249\u00d74=\"; New = \"733\u00d74=\" },\n    @{ Old = \"477\u00d76=\"; New = \"853\u00d77=\" },\n    @{ Old = \"788\u00d75=\"; New = \"871\u00d79=\" },\n    @{ Old = \"988\u00d75=\"; New = \"321\u00d74=\" },\n    @{ Old = \"546\u00d76=\"; New = \"553\u00d73=\" },\n    @{ Old = \"120\u00d73=\"; New = \"803\u00d75=\" },\n    @{ Old = \"551\u00d72=\"; New = \"965\u00d76=\" },\n    @{ Old = \"634\u00d75=\"; New = \"361\u00d73=\" },\n    @{ Old = \"996\u00d75=\"; New = \"126\u00d78=\" },\n    @{ Old = \"648\u00d79=\"; New = \"996\u00d74=\" },\n    @{ Old = \"631\u00d74=\"; New = \"418\u00d74=\" },\n    @{ Old = \"597\u00d77=\"; New = \"710\u00d76=\" },\n    @{ Old = \"974\u00d77=\"; New = \"166\u00d79=\" },\n    @{ Old = \"580\u00d78=\"; New = \"486\u00d74=\" },\n    @{ Old = \"175\u00d72=\"; New = \"654\u00d75=\" },\n    @{ Old = \"250\u00d76=\"; New = \"538\u00d76=\" },\n    @{ Old = \"576\u00d78=\"; New = \"486\u00d79=\" },\n    @{ Old = \"141\u00d73=\"; New = \"695\u00d75=\" },\n    @{ Old = \"390\u00d75=\"; New = \"147\u00d76=\" },\n    @{ Old = \"863\u00d78=\"; New = \"611\u00d73=\" },\n    @{ Old = \"672\u00d72=\"; New = \"678\u00d72=\" },\n    @{ Old = \"813\u00d75=\"; New = \"435\u00d76=\" },\n    @{ Old = \"415\u00d78=\"; New = \"483\u00d77=\" },\n    @{ Old = \"378\u00d77=\"; New = \"186\u00d73=\" }\n)\n\nforeach ($r in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $r.Old\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $r.New\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.MatchSoundsLike = $false\n    $find.MatchAllWordForms = $false\n    $find.Execute($r.Old, $false, $false, $false, $false, $false, $true, 1, $false, $r.New, 2)\n}\n\nWrite-Output \"done\"\n"}
